# live_trading_results.xlsx - record Trade #13 closing and Trade #31 opening.
#
# Trade #13 (leadlag row 12 / newly appended All Trades row 14) is closed
# with an exit price / pnl, and a brand-new Trade #31 is appended to the
# leadlag sheet as a freshly-opened (still-OPEN) trade. The Summary and
# Comparison roll-up sheets are refreshed to reflect the new trade counts
# and aggregate stats.

$wb = $excel.ActiveWorkbook

# Helper: force a cell to hold a literal text value (percent-looking
# strings like "61.5%" and date-looking strings like "2026-02-16" would
# otherwise get auto-coerced into numbers/dates by Excel's input parser).
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 13                                 # C2 Total Trades (OVERALL)
Set-TextValue $summary.Cells.Item(2, 4) "61.5%"                      # D2 Win Rate
Set-TextValue $summary.Cells.Item(2, 5) "+2.2450%"                   # E2 Total P&L %
Set-TextValue $summary.Cells.Item(2, 6) "+0.1727%"                   # F2 Avg Trade

$summary.Cells.Item(3, 3).Value = 25                                 # C3 Total Trades (leadlag)
Set-TextValue $summary.Cells.Item(3, 4) "28.0%"                      # D3 Win Rate
Set-TextValue $summary.Cells.Item(3, 5) "+2.1997%"                   # E3 Total P&L %
Set-TextValue $summary.Cells.Item(3, 6) "+0.0880%"                   # F3 Avg Trade

# ---------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Trade #13 (row 12) transitions from OPEN to CLOSED.
$leadlag.Cells.Item(12, 7).Value = 69028.823174                      # G12 Exit Price
Set-TextValue $leadlag.Cells.Item(12, 8) "CLOSED"                    # H12 Status
$leadlag.Cells.Item(12, 9).Value = -0.5831                           # I12 P&L %
$leadlag.Cells.Item(12, 10).Value = -5.83                            # J12 P&L $
Set-TextValue $leadlag.Cells.Item(12, 13) "time_exit_5min"           # M12 Exit Reason
$leadlag.Cells.Item(12, 14).Value = 5                                # N12 Duration (min)

# New row 27: Trade #31, freshly opened.
$leadlag.Cells.Item(27, 1).Value = 31                                 # A27 Trade #
Set-TextValue $leadlag.Cells.Item(27, 2) "2026-02-16"                 # B27 Date
Set-TextValue $leadlag.Cells.Item(27, 3) "21:28:05"                   # C27 Time
Set-TextValue $leadlag.Cells.Item(27, 4) "leadlag"                    # D27 Strategy
Set-TextValue $leadlag.Cells.Item(27, 5) "UP"                         # E27 Side
$leadlag.Cells.Item(27, 6).Value = 68945.64                           # F27 Entry Price
$leadlag.Cells.Item(27, 7).Value = ""                                 # G27 Exit Price (blank, still open)
Set-TextValue $leadlag.Cells.Item(27, 8) "OPEN"                       # H27 Status
$leadlag.Cells.Item(27, 9).Value = 0                                  # I27 P&L %
$leadlag.Cells.Item(27, 10).Value = 0                                 # J27 P&L $
$leadlag.Cells.Item(27, 11).Value = 0.75                              # K27 Confidence
Set-TextValue $leadlag.Cells.Item(27, 12) "Coinbase leading with 0.088% move"  # L27 Entry Reason
$leadlag.Cells.Item(27, 13).Value = ""                                # M27 Exit Reason (blank, still open)
$leadlag.Cells.Item(27, 14).Value = 0                                 # N27 Duration (min)

# ---------------------------------------------------------------------
# All Trades sheet - append the newly-closed Trade #13.
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(14, 1).Value = 13                                # A14 Trade #
Set-TextValue $allTrades.Cells.Item(14, 2) "2026-02-16"                # B14 Date
Set-TextValue $allTrades.Cells.Item(14, 3) "21:22:51"                  # C14 Time
Set-TextValue $allTrades.Cells.Item(14, 4) "leadlag"                   # D14 Strategy
Set-TextValue $allTrades.Cells.Item(14, 5) "UP"                        # E14 Side
$allTrades.Cells.Item(14, 6).Value = 69433.7                           # F14 Entry Price
$allTrades.Cells.Item(14, 7).Value = 69028.823174                      # G14 Exit Price
Set-TextValue $allTrades.Cells.Item(14, 8) "CLOSED"                    # H14 Status
$allTrades.Cells.Item(14, 9).Value = -0.5831                           # I14 P&L %
$allTrades.Cells.Item(14, 10).Value = -5.83                            # J14 P&L $
$allTrades.Cells.Item(14, 11).Value = 0.75                             # K14 Confidence
Set-TextValue $allTrades.Cells.Item(14, 12) "Binance leading with 0.078% move"  # L14 Entry Reason
Set-TextValue $allTrades.Cells.Item(14, 13) "time_exit_5min"           # M14 Exit Reason
$allTrades.Cells.Item(14, 14).Value = 5                                # N14 Duration (min)

# ---------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(2, 2).Value = 25                                # B2 Total Trades
Set-TextValue $comparison.Cells.Item(2, 3) "28.0%"                     # C2 Win Rate
Set-TextValue $comparison.Cells.Item(2, 4) "2.89"                      # D2 Profit Factor
Set-TextValue $comparison.Cells.Item(2, 6) "-0.2915%"                  # F2 Avg Loss %
Set-TextValue $comparison.Cells.Item(2, 7) "1.65"                      # G2 Win/Loss Ratio
Set-TextValue $comparison.Cells.Item(2, 8) "-0.5831%"                  # H2 Max Drawdown
